# Update "last_edited_time" (column D) for rows 7-12 and the numeric
# property values in row 7 (Tháng 7), per commit:
#   "xoa cac phan tong cua ti le chiet khau"
# (the discount-rate related totals for July were recalculated/updated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2024-07-25T15:02:00.000Z"

# Update last_edited_time for rows 7 through 12
foreach ($r in 7..12) {
    $ws.Cells.Item($r, 4).Value = $newTimestamp
}

# Update numeric property values for row 7 (Tháng 7)
$ws.Range("W7").Value = 269931000
$ws.Range("AA7").Value = 277037000
$ws.Range("AE7").Value = 546968000
$ws.Range("AH7").Value = 468968000
$ws.Range("AK7").Value = 71
$ws.Range("AN7").Value = 78000000
$ws.Range("AQ7").Value = 525268000
